$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: update title and link
$ws.Range("D4").Value = "뽀로로(PORORO) 자연어처리 라이브러리 활용기 (dev. and maintained by 카카오 브레인(Kakao Brain) PORORO팀)"
$ws.Range("E4").Value = "https://teddylee777.github.io/machine-learning/nlp-korean-pororo"

# Row 29: update title
$ws.Range("D29").Value = "프로메디우스"

# Row 37: update title and link
$ws.Range("D37").Value = "[Paper Review] Few-Shot Anomaly Detection"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1874&mod=document&pageid=1"

# Row 51: update title and link
$ws.Range("D51").Value = "[윈도우10] 윈도우 화면잠금 단축키는 무엇일까?"
$ws.Range("E51").Value = "https://bskyvision.com/1231"
